$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K24").Value = [double]"-1.284357670755455"
$ws.Range("J25").Value = [double]"-0.1748886306845208"
$ws.Range("K25").Value = [double]"1.091677354385908"
$ws.Range("I26").Value = [double]"-0.560290197271857"
$ws.Range("J26").Value = [double]"0.7062757877985717"
$ws.Range("K26").Value = [double]"8.545046978308823E-06"
$ws.Range("H27").Value = [double]"0.8677816465512651"
$ws.Range("I27").Value = [double]"2.134347631621694"
$ws.Range("J27").Value = [double]"1.4280803888701"
$ws.Range("K27").Value = [double]"-0.181580002546667"
$ws.Range("G28").Value = [double]"-1.060290197271857"
$ws.Range("H28").Value = [double]"0.2062757877985718"
$ws.Range("I28").Value = [double]"-0.4999914549530217"
$ws.Range("J28").Value = [double]"-2.109651846369789"
$ws.Range("K28").Value = [double]"1.918509177329452"
$ws.Range("F29").Value = [double]"-0.2579984325430758"
$ws.Range("G29").Value = [double]"1.008567552527353"
$ws.Range("H29").Value = [double]"0.3023003097757595"
$ws.Range("I29").Value = [double]"-1.307360081641008"
$ws.Range("J29").Value = [double]"2.720800942058233"
$ws.Range("K29").Value = [double]"-0.4557938899696836"
$ws.Range("E30").Value = [double]"-0.005406304358785974"
$ws.Range("F30").Value = [double]"1.261159680711643"
$ws.Range("G30").Value = [double]"0.5548924379600493"
$ws.Range("H30").Value = [double]"-1.054767953456718"
$ws.Range("I30").Value = [double]"2.973393070242523"
$ws.Range("J30").Value = [double]"-0.2032017617853938"
$ws.Range("K30").Value = [double]"0.6864325751146225"
$ws.Range("D31").Value = [double]"-0.4108264758588121"
$ws.Range("E31").Value = [double]"0.8557395092116167"
$ws.Range("F31").Value = [double]"0.1494722664600232"
$ws.Range("G31").Value = [double]"-1.460188124956744"
$ws.Range("H31").Value = [double]"2.567972898742497"
$ws.Range("I31").Value = [double]"-0.6086219332854199"
$ws.Range("J31").Value = [double]"0.2810124036145964"
$ws.Range("K31").Value = [double]"-1.994615918782514"
$ws.Range("C32").Value = [double]"-0.560290197271857"
$ws.Range("D32").Value = [double]"0.7062757877985717"
$ws.Range("E32").Value = [double]"8.545046978308823E-06"
$ws.Range("F32").Value = [double]"-1.609651846369789"
$ws.Range("G32").Value = [double]"2.418509177329452"
$ws.Range("H32").Value = [double]"-0.7580856546984649"
$ws.Range("I32").Value = [double]"0.1315486822015515"
$ws.Range("J32").Value = [double]"-2.144079640195559"
$ws.Range("K32").Value = [double]"0.2000055876513443"
$ws.Range("B33").Value = [double]"0.6898948471401811"
$ws.Range("C33").Value = [double]"1.95646083221061"
$ws.Range("D33").Value = [double]"1.250193589459016"
$ws.Range("E33").Value = [double]"-0.3594668019577509"
$ws.Range("F33").Value = [double]"3.66869422174149"
$ws.Range("G33").Value = [double]"0.4920993897135733"
$ws.Range("H33").Value = [double]"1.38173372661359"
$ws.Range("I33").Value = [double]"-0.8938945957835209"
$ws.Range("J33").Value = [double]"1.450190632063382"
$ws.Range("K33").Value = [double]"0.8268227261614813"
$ws.Range("B34").Value = [double]"1.266565985070429"
$ws.Range("C34").Value = [double]"0.5602987423188353"
$ws.Range("D34").Value = [double]"-1.049361649097932"
$ws.Range("E34").Value = [double]"2.978799374601309"
$ws.Range("F34").Value = [double]"-0.1977954574266079"
$ws.Range("G34").Value = [double]"0.6918388794734085"
$ws.Range("H34").Value = [double]"-1.583789442923702"
$ws.Range("I34").Value = [double]"0.7602957849232013"
$ws.Range("J34").Value = [double]"0.1369278790213002"
$ws.Range("K34").Value = [double]"5.647645974903085"
$ws.Range("B35").Value = [double]"-0.7062672427515935"
$ws.Range("C35").Value = [double]"-2.315927634168361"
$ws.Range("D35").Value = [double]"1.71223338953088"
$ws.Range("E35").Value = [double]"-1.464361442497037"
$ws.Range("F35").Value = [double]"-0.5747271055970202"
$ws.Range("G35").Value = [double]"-2.850355427994131"
$ws.Range("H35").Value = [double]"-0.5062702001472275"
$ws.Range("I35").Value = [double]"-1.129638106049129"
$ws.Range("J35").Value = [double]"4.381079989832656"
$ws.Range("K35").Value = [double]"-0.5837405679867658"
$ws.Range("B36").Value = [double]"-1.609660391416767"
$ws.Range("C36").Value = [double]"2.418500632282474"
$ws.Range("D36").Value = [double]"-0.7580941997454431"
$ws.Range("E36").Value = [double]"0.1315401371545732"
$ws.Range("F36").Value = [double]"-2.144088185242537"
$ws.Range("G36").Value = [double]"0.199997042604366"
$ws.Range("H36").Value = [double]"-0.423370863297535"
$ws.Range("I36").Value = [double]"5.08734723258425"
$ws.Range("J36").Value = [double]"0.1225266747648277"
$ws.Range("K36").Value = [double]"-0.0232082318966747"
$ws.Range("B37").Value = [double]"4.028161023699241"
$ws.Range("C37").Value = [double]"0.8515661916713242"
$ws.Range("D37").Value = [double]"1.741200528571341"
$ws.Range("E37").Value = [double]"-0.53442779382577"
$ws.Range("F37").Value = [double]"1.809657434021133"
$ws.Range("G37").Value = [double]"1.186289528119232"
$ws.Range("H37").Value = [double]"6.697007624001017"
$ws.Range("I37").Value = [double]"1.732187066181595"
$ws.Range("J37").Value = [double]"1.586452159520093"
$ws.Range("K37").Value = [double]"1.462338942638652"
$ws.Range("B38").Value = [double]"-3.176594832027917"
$ws.Range("C38").Value = [double]"-2.286960495127901"
$ws.Range("D38").Value = [double]"-4.562588817525011"
$ws.Range("E38").Value = [double]"-2.218503589678108"
$ws.Range("F38").Value = [double]"-2.841871495580009"
$ws.Range("G38").Value = [double]"2.668846600301776"
$ws.Range("H38").Value = [double]"-2.295973957517646"
$ws.Range("I38").Value = [double]"-2.441708864179148"
$ws.Range("J38").Value = [double]"-2.565822081060589"
$ws.Range("K38").Value = [double]"-1.791630897729207"
$ws.Range("B39").Value = [double]"0.8896343369000164"
$ws.Range("C39").Value = [double]"-1.385993985497094"
$ws.Range("D39").Value = [double]"0.9580912423498091"
$ws.Range("E39").Value = [double]"0.3347233364479081"
$ws.Range("F39").Value = [double]"5.845441432329693"
$ws.Range("G39").Value = [double]"0.8806208745102708"
$ws.Range("H39").Value = [double]"0.7348859678487685"
$ws.Range("I39").Value = [double]"0.6107727509673282"
$ws.Range("J39").Value = [double]"1.38496393429871"
$ws.Range("K39").Value = [double]"-0.07658964140479418"
$ws.Range("B40").Value = [double]"-2.275628322397111"
$ws.Range("C40").Value = [double]"0.06845690544979277"
$ws.Range("D40").Value = [double]"-0.5549110004521083"
$ws.Range("E40").Value = [double]"4.955807095429677"
$ws.Range("F40").Value = [double]"-0.009013462389745541"
$ws.Range("G40").Value = [double]"-0.1547483690512479"
$ws.Range("H40").Value = [double]"-0.2788615859326882"
$ws.Range("I40").Value = [double]"0.495329597398694"
$ws.Range("J40").Value = [double]"-0.9662239783048105"
$ws.Range("K40").Value = [double]"-0.3659156202749045"
$ws.Range("B41").Value = [double]"2.344085227846903"
$ws.Range("C41").Value = [double]"1.720717321945002"
$ws.Range("D41").Value = [double]"7.231435417826788"
$ws.Range("E41").Value = [double]"2.266614860007365"
$ws.Range("F41").Value = [double]"2.120879953345863"
$ws.Range("G41").Value = [double]"1.996766736464422"
$ws.Range("H41").Value = [double]"2.770957919795805"
$ws.Range("I41").Value = [double]"1.3094043440923"
$ws.Range("J41").Value = [double]"1.909712702122206"
$ws.Range("K41").Value = [double]"1.928070503448779"
$ws.Range("B42").Value = [double]"-0.6233679059019011"
$ws.Range("C42").Value = [double]"4.887350189979884"
$ws.Range("D42").Value = [double]"-0.07747036783953831"
$ws.Range("E42").Value = [double]"-0.2232052745010407"
$ws.Range("F42").Value = [double]"-0.347318491382481"
$ws.Range("G42").Value = [double]"0.4268726919489012"
$ws.Range("H42").Value = [double]"-1.034680883754603"
$ws.Range("I42").Value = [double]"-0.4343725257246973"
$ws.Range("J42").Value = [double]"-0.4160147243981243"
$ws.Range("K42").Value = [double]"0.3423398389707444"
$ws.Range("B43").Value = [double]"5.510718095881785"
$ws.Range("C43").Value = [double]"0.5458975380623627"
$ws.Range("D43").Value = [double]"0.4001626314008604"
$ws.Range("E43").Value = [double]"0.2760494145194201"
$ws.Range("F43").Value = [double]"1.050240597850802"
$ws.Range("G43").Value = [double]"-0.4113129778527023"
$ws.Range("H43").Value = [double]"0.1889953801772037"
$ws.Range("I43").Value = [double]"0.2073531815037768"
$ws.Range("J43").Value = [double]"0.9657077448726454"
$ws.Range("K43").Value = [double]"-0.1685777075333164"
$ws.Range("B44").Value = [double]"-4.964820557819422"
$ws.Range("C44").Value = [double]"-5.110555464480925"
$ws.Range("D44").Value = [double]"-5.234668681362365"
$ws.Range("E44").Value = [double]"-4.460477498030983"
$ws.Range("F44").Value = [double]"-5.922031073734487"
$ws.Range("G44").Value = [double]"-5.321722715704581"
$ws.Range("H44").Value = [double]"-5.303364914378008"
$ws.Range("I44").Value = [double]"-4.54501035100914"
$ws.Range("J44").Value = [double]"-5.679295803415101"
$ws.Range("B45").Value = [double]"-0.1457349066615024"
$ws.Range("C45").Value = [double]"-0.2698481235429427"
$ws.Range("D45").Value = [double]"0.5043430597884395"
$ws.Range("E45").Value = [double]"-0.957210515915065"
$ws.Range("F45").Value = [double]"-0.356902157885159"
$ws.Range("G45").Value = [double]"-0.3385443565585859"
$ws.Range("H45").Value = [double]"0.4198102068102827"
$ws.Range("I45").Value = [double]"-0.7144752455956791"
$ws.Range("B46").Value = [double]"-0.1241132168814403"
$ws.Range("C46").Value = [double]"0.6500779664499419"
$ws.Range("D46").Value = [double]"-0.8114756092535627"
$ws.Range("E46").Value = [double]"-0.2111672512236566"
$ws.Range("F46").Value = [double]"-0.1928094498970836"
$ws.Range("G46").Value = [double]"0.5655451134717852"
$ws.Range("H46").Value = [double]"-0.5687403389341767"
$ws.Range("B47").Value = [double]"0.7741911833313821"
$ws.Range("C47").Value = [double]"-0.6873623923721224"
$ws.Range("D47").Value = [double]"-0.08705403434221631"
$ws.Range("E47").Value = [double]"-0.06869623301564326"
$ws.Range("F47").Value = [double]"0.6896583303532254"
$ws.Range("G47").Value = [double]"-0.4446271220527364"
$ws.Range("B48").Value = [double]"-1.461553575703505"
$ws.Range("C48").Value = [double]"-0.8612452176735985"
$ws.Range("D48").Value = [double]"-0.8428874163470255"
$ws.Range("E48").Value = [double]"-0.08453285297815677"
$ws.Range("F48").Value = [double]"-1.218818305384119"
$ws.Range("B49").Value = [double]"0.600308358029906"
$ws.Range("C49").Value = [double]"0.6186661593564791"
$ws.Range("D49").Value = [double]"1.377020722725348"
$ws.Range("E49").Value = [double]"0.2427352703193859"
$ws.Range("B50").Value = [double]"0.01835780132657305"
$ws.Range("C50").Value = [double]"0.7767123646954417"
$ws.Range("D50").Value = [double]"-0.3575730877105201"
$ws.Range("B51").Value = [double]"0.7583545633688686"
$ws.Range("C51").Value = [double]"-0.3759308890370932"
$ws.Range("B52").Value = [double]"-1.134285452405962"
